$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) target cells to Text format first so that
# numeric-looking strings (e.g. "1.007") are stored as literal text
# instead of being parsed into numbers by the input parser.
$priceCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cell values.
$ws.Range("D2").Value = "28.162.56"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "1.877.63"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "313.46"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "0.5126"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "0.3905"
$ws.Range("E8").Value = "  +2.80%  "
$ws.Range("D9").Value = "0.08325"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").Value = "1.120"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").Value = "41.46"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "6.214"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").Value = "20.66"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "1.867.06"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "7.264"
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "0.00001101"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "91.09"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "0.06647"
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").Value = "17.76"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "6.031"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").Value = "28.177.66"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "11.12"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("D25").Value = "2.271"
$ws.Range("E25").Value = "  +2.05%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.085.18"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.499"
$ws.Range("E27").Value = "  -2.75%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "158.77"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "20.60"
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "125.10"
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.1064"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.040"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.863"
$ws.Range("E33").Value = "  +5.13%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.595"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "9.690"
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "0.02455"
$ws.Range("E36").Value = "  +2.29%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.06533"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.2186"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "1.204"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6502"
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.233"
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("B42").Value = "InternetComputer(DFINITY)"
$ws.Range("C42").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D42").Value = "4.988"
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "11.27"
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "0.6134"
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.11"
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "1.285"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.667"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "2.008"
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "1.233"
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "121.22"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "78.57"
$ws.Range("E51").Value = "  -1.47%  "
